$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.611
$ws.Range("D3").Value = -8.005000000000001
$ws.Range("D14").Value = -8.005000000000001
$ws.Range("D21").Value = -8.1
$ws.Range("D23").Value = -7.874
$ws.Range("D25").Value = -7.842999999999999
$ws.Range("E25").Value = 17.104
$ws.Range("D26").Value = -7.635
$ws.Range("E27").Value = 16.764
$ws.Range("D29").Value = -7.325
$ws.Range("E31").Value = 17.257
$ws.Range("E39").Value = 16.573
$ws.Range("E48").Value = 17.179
$ws.Range("E51").Value = 16.617
$ws.Range("E52").Value = 16.543
$ws.Range("D53").Value = -7.755
$ws.Range("E55").Value = 16.416
$ws.Range("E56").Value = 16.276
$ws.Range("D57").Value = -7.923999999999999
$ws.Range("E57").Value = 16.568
$ws.Range("D59").Value = -8.061
$ws.Range("D69").Value = -7.597
$ws.Range("E73").Value = 16.572
$ws.Range("D79").Value = -7.885
$ws.Range("D83").Value = -8.237
$ws.Range("E89").Value = 17.57
$ws.Range("E90").Value = 16.73
$ws.Range("D91").Value = -6.87
$ws.Range("E92").Value = 17.253
$ws.Range("D93").Value = -7.595000000000001
